# This script re-generates the quadratic/linear problem data tables
# ("volver a generar problemas cuadraticos y lineales") by writing new
# text values into the cells that hold expressions / numeric results.
#
# Many of these cells look like numbers (e.g. "0.77", "-104.6") but are
# actually stored as *text* in the workbook (so that values such as
# "6.8999999999999995" keep their exact textual representation instead
# of being rounded/normalized as IEEE-754 doubles). Assigning a plain
# numeric-looking string via .Value would make Excel auto-convert it to
# a real number, so we briefly force the cell to Text format, assign the
# value, and then clear the cell formatting again (this mirrors how the
# original cells have no explicit style, just General format, while
# still keeping the value as text).
#
# Sheets are addressed by their (1-based) tab index rather than by name:
# this workbook has both "Vector_bf" and "Vector_BF" as sheet names, and
# name-based lookups are case-insensitive, so `Worksheets.Item("Vector_BF")`
# would otherwise collide with `Worksheets.Item("Vector_bf")`.

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 3: Restricciones_del_follower
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

Set-TextValue $ws.Range("A2") "1.7500000000000169 - 2x_1 + y_1 - y_2"
Set-TextValue $ws.Range("B2") "0.7499999999999831"
Set-TextValue $ws.Range("D2") "0.77"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "6.1"

Set-TextValue $ws.Range("A3") "6.549999999999994 + x_1 - 3x_2 + y_2"
Set-TextValue $ws.Range("B3") "-8.549999999999994"
Set-TextValue $ws.Range("D3") "0.46"
Set-TextValue $ws.Range("E3") "6.8999999999999995"

Set-TextValue $ws.Range("A4") "104.6 - y_1"
Set-TextValue $ws.Range("B4") "-104.6"
Set-TextValue $ws.Range("D4") "0.41"
Set-TextValue $ws.Range("E4") "-5.5"
Set-TextValue $ws.Range("F4") "-5.1"

Set-TextValue $ws.Range("A5") "-2.05 - y_2"
Set-TextValue $ws.Range("B5") "-2.05"
Set-TextValue $ws.Range("D5") "0.64"
Set-TextValue $ws.Range("E5") "-9.8"
Set-TextValue $ws.Range("F5") "-8.5"

# ---------------------------------------------------------------
# 4: Punto_modificado
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

Set-TextValue $ws.Range("A2") "52.150000000000006"
Set-TextValue $ws.Range("B2") "20.25"
Set-TextValue $ws.Range("C2") "104.6"
Set-TextValue $ws.Range("D2") "2.05"

# ---------------------------------------------------------------
# 5: Vector_bf
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

Set-TextValue $ws.Range("A2") "3.64"
Set-TextValue $ws.Range("A3") "-0.050000000000000044"

# ---------------------------------------------------------------
# 6: Vector_BF
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

Set-TextValue $ws.Range("A2") "-4.8999999999999995"
Set-TextValue $ws.Range("A3") "19.7"
Set-TextValue $ws.Range("A4") "-6.0"
Set-TextValue $ws.Range("A5") "-16.7"
